# recommender.xlsx edit: move dataset and models
# Insert 4 new columns (B:E) for "data", "LR", "gamma", "epoch" before the
# existing "milestone" column, add two new style number-formats, backfill
# the new parameter columns for the existing rows, and append five new
# benchmark rows (9-13) plus their hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Clear the existing hyperlinks up front - their target cells are about
#    to shift right by four columns, and it is simplest to re-create them
#    once the new layout is final.
# ---------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 2. Insert four new columns at B:E - this shifts the old B (milestone),
#    C (Loss) and D:G (Acc@1/5/10/20) columns to F, G and H:K.
# ---------------------------------------------------------------------
$ws.Range("B1:E1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 3. New header row.
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "data"
$ws.Range("C1").Value = "LR"
$ws.Range("D1").Value = "gamma"
$ws.Range("E1").Value = "epoch"

# ---------------------------------------------------------------------
# 4. Backfill the new parameter columns for the pre-existing rows
#    (2-7), which all share the same data/gamma/epoch values.
# ---------------------------------------------------------------------
$rows27 = @(2, 3, 4, 5, 6, 7)
foreach ($r in $rows27) {
  $ws.Cells.Item($r, 2).Value = 5
  $ws.Cells.Item($r, 5).Value = 1000
  $ws.Cells.Item($r, 4).Value = 0.1
}
$ws.Range("C2").Value = 0.001
$ws.Range("C3").Value = 0.01
$ws.Range("C4").Value = 0.1
$ws.Range("C5").Value = 0.01
$ws.Range("C6").Value = 0.01
$ws.Range("C7").Value = 0.01

# ---------------------------------------------------------------------
# 5. New rows 9-13.
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "EP200"
$ws.Range("B9").Value = 2293
$ws.Range("C9").Value = 0.001
$ws.Range("D9").Value = 0.1
$ws.Range("E9").Value = 200
$ws.Range("F9").Value = "(200,)"
$ws.Range("G9").Value = 0.1054
$ws.Range("H9").Value = 0.037069
$ws.Range("I9").Value = 0.089984
$ws.Range("J9").Value = 0.124146
$ws.Range("K9").Value = 0.165286

$ws.Range("A10").Value = "EP20"
$ws.Range("B10").Value = 2293
$ws.Range("C10").Value = 0.01
$ws.Range("D10").Value = 0.1
$ws.Range("E10").Value = 20
$ws.Range("F10").Value = "(20,)"
$ws.Range("G10").Value = 0.141675
$ws.Range("H10").Value = 0.018026
$ws.Range("I10").Value = 0.040558
$ws.Range("J10").Value = 0.057276
$ws.Range("K10").Value = 0.075011

$ws.Range("A11").Value = "LR0005"
$ws.Range("B11").Value = 2293
$ws.Range("C11").Value = 0.005
$ws.Range("D11").Value = 0.1
$ws.Range("E11").Value = 20
$ws.Range("F11").Value = "(20,)"
$ws.Range("G11").Value = 0.125113
$ws.Range("H11").Value = 0.023695
$ws.Range("I11").Value = 0.06091
$ws.Range("J11").Value = 0.081553
$ws.Range("K11").Value = 0.106265

$ws.Range("A12").Value = "LR0002"
$ws.Range("B12").Value = 2293
$ws.Range("C12").Value = 0.002
$ws.Range("D12").Value = 0.1
$ws.Range("E12").Value = 20
$ws.Range("F12").Value = "(20,)"
$ws.Range("G12").Value = 0.111472
$ws.Range("H12").Value = 0.031981
$ws.Range("I12").Value = 0.076755
$ws.Range("J12").Value = 0.105539
$ws.Range("K12").Value = 0.138247

$ws.Range("A13").Value = "EP150"
$ws.Range("B13").Value = 2293
$ws.Range("C13").Value = 0.001
$ws.Range("D13").Value = 0.5
$ws.Range("E13").Value = 150
$ws.Range("F13").Value = "(10,)"
$ws.Range("G13").Value = 0.101838
$ws.Range("H13").Value = 0.041867
$ws.Range("I13").Value = 0.103649
$ws.Range("J13").Value = 0.136648
$ws.Range("K13").Value = 0.17866

# ---------------------------------------------------------------------
# 6. Styling: columns B:E (new numeric params) get the plain
#    horizontal-left style; column G (Loss) gets a brand-new
#    accounting-style number format (red negatives); columns H:K keep
#    the original Acc@ number format.
# ---------------------------------------------------------------------
$ws.Range("B2:E13").HorizontalAlignment = -4131

$ws.Range("G1:G13").NumberFormat = "0.000000_);[Red]\(0.000000\)"
$ws.Range("G1:G13").HorizontalAlignment = -4131

$ws.Range("H7:K7").NumberFormat = "0.000000_ "
$ws.Range("H7:K13").NumberFormat = "0.000000_ "
$ws.Range("H7:K13").HorizontalAlignment = -4131

# ---------------------------------------------------------------------
# 7. Column widths to roughly match the new layout.
# ---------------------------------------------------------------------
$ws.Columns("B:E").ColumnWidth = 10.75
$ws.Columns("F:F").ColumnWidth = 10.75

# ---------------------------------------------------------------------
# 8. Re-create the hyperlinks on the (now shifted) header cells, plus the
#    new one on J11.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("H1"), "mailto:Acc@1", "", "")
$ws.Hyperlinks.Add($ws.Range("I1:K1"), "mailto:Acc@1", "", "", "Acc@1")
$ws.Hyperlinks.Add($ws.Range("I1"), "mailto:Acc@5", "", "")
$ws.Hyperlinks.Add($ws.Range("J1"), "mailto:Acc@10", "", "")
$ws.Hyperlinks.Add($ws.Range("K1"), "mailto:Acc@20", "", "")
$ws.Hyperlinks.Add($ws.Range("J11"), "mailto:Acc@10", "", "Acc@10: 0.081553")

Write-Host "edit complete"
